$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "Max_land_usage" worksheet as the LAST sheet in the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Max_land_usage"

# ---------------------------------------------------------------------------
# Row 1 — grouping headers (bold, bordered, centered, merged)
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Tech_category"
$ws.Cells.Item(1, 2).Value = "Supply"
$ws.Cells.Item(1, 9).Value = "Conversion"
$ws.Cells.Item(1, 12).Value = "Transmission"

$ws.Range("B1:H1").Merge()
$ws.Range("I1:K1").Merge()

# ---------------------------------------------------------------------------
# Row 2 — technology headers
# ---------------------------------------------------------------------------
$row2 = @("Technology","Natural_gas_supply","Oil_supply","BW_supply","PV_PP","Wind_PP","Geo_PP","Hydro_PP","HFO_PP","OCGT_PP","BW_PP","Elec_transmission_distribution")
for ($c = 1; $c -le $row2.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
}

# ---------------------------------------------------------------------------
# Row 3 — "Years" label
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Years"

# ---------------------------------------------------------------------------
# Rows 4-14 — Y0..Y10 labels in column A, data values 1E+20 in columns B:L
# ---------------------------------------------------------------------------
$years = @("Y0","Y1","Y2","Y3","Y4","Y5","Y6","Y7","Y8","Y9","Y10")
$bigVal = [double]"1E+20"
for ($r = 0; $r -lt $years.Length; $r++) {
    $rowIdx = 4 + $r
    $ws.Cells.Item($rowIdx, 1).Value = $years[$r]
    for ($c = 2; $c -le 12; $c++) {
        $ws.Cells.Item($rowIdx, $c).Value = $bigVal
    }
}

# ---------------------------------------------------------------------------
# Formatting — bold font + thin box border + center/top alignment on the
# header block (rows 1-2, columns A:L), the "Years" label (A3) and on
# column A for the Y0..Y10 data rows (A4:A14).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:L2")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$yearsLabel = $ws.Range("A3")
$yearsLabel.Font.Bold = $true
$yearsLabel.Borders.LineStyle = 1
$yearsLabel.HorizontalAlignment = -4108
$yearsLabel.VerticalAlignment = -4160

$labelColumn = $ws.Range("A4:A14")
$labelColumn.Font.Bold = $true
$labelColumn.Borders.LineStyle = 1
$labelColumn.HorizontalAlignment = -4108
$labelColumn.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Page margins — match the workbook-standard margins used by every other
# parameter sheet (0.75in/0.75in/1in/1in/0.5in/0.5in).
# ---------------------------------------------------------------------------
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Make the new sheet the active/selected sheet (matches tabSelected moving
# from capacity_factor_resource onto the newly-created final sheet).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select()
